$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the values currently in row 2 (the latest price entry, 04-11-2025)
# before inserting a brand-new row above it.
$desc     = $ws.Range("B2").Value2
$grade    = $ws.Range("C2").Value2
$price    = $ws.Range("D2").Value2
$circDate = $ws.Range("E2").Value2
$circLink = $ws.Range("F2").Value2

# Insert a new row at position 2; existing rows (2..147) shift down to (3..148),
# carrying their values, styles and hyperlinks with them.
$ws.Rows("2:2").Insert()

# The freshly inserted row picks up formatting from the header row above it,
# so first re-apply the plain data-row formatting (from row 3, the shifted
# former row 2) across the whole row.
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A2:F2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill the newly inserted row 2 with the latest entry (05-11-2025), reusing
# the description/grade/price/circular date/link from the prior latest entry
# (the underlying PDF circular itself was not reissued that day).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "05-11-2025"

$ws.Range("B2").Value = $desc
$ws.Range("C2").Value = $grade
$ws.Range("D2").Value = $price

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = $circDate

$ws.Range("F2").Value = $circLink

# Forcing "Text" number format on A2/E2 (so the date-like strings aren't
# reinterpreted as dates) mints a new style; copy the original style back
# on top now that the literal text values are safely in place.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wb.Save()
